# Update of league bases, from 13-06-2024 19:35.
# Two pairs of adjacent match rows got swapped in the
# "Kazakhstan Premier League" sheet: rows 88<->89 and rows 148<->149.
# Column A (the running row index) is left untouched; every other
# populated cell (columns B, and E through AD) of each pair of rows is
# exchanged with its counterpart - the id in column B simply moves to
# the other row along with the rest of that match's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B..AD are spreadsheet columns 2..30 (column C/D are identical
# between the paired rows, so swapping them is harmless/idempotent too).
$firstCol = 2
$lastCol = 30

$pairs = @(@(88, 89), @(148, 149))

foreach ($pair in $pairs) {
    $row1 = $pair[0]
    $row2 = $pair[1]

    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $cell1 = $ws.Cells.Item($row1, $col)
        $cell2 = $ws.Cells.Item($row2, $col)

        $val1 = $cell1.Value2
        $val2 = $cell2.Value2

        $cell1.Value = $val2
        $cell2.Value = $val1
    }
}
